$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.239.25"
$ws.Range("E2").Value = "  -3.43%  "
$ws.Range("D3").Value = "3.177.69"
$ws.Range("E3").Value = "  -8.39%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'564.88"
$ws.Range("E5").Value = "  -3.62%  "
$ws.Range("D6").Value = "'168.54"
$ws.Range("E6").Value = "  -5.22%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.605"
$ws.Range("E8").Value = "  -3.54%  "
$ws.Range("D9").Value = "3.176.39"
$ws.Range("E9").Value = "  -8.33%  "
$ws.Range("E10").Value = "  -6.90%  "
$ws.Range("E11").Value = "  -5.38%  "
$ws.Range("D12").Value = "'0.396"
$ws.Range("E12").Value = "  -5.56%  "
$ws.Range("D13").Value = "3.721.13"
$ws.Range("E13").Value = "  -8.56%  "
$ws.Range("E14").Value = "  +1.42%  "
$ws.Range("D15").Value = "'27.34"
$ws.Range("E15").Value = "  -9.60%  "
$ws.Range("D16").Value = "64.229.19"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("E17").Value = "  -5.51%  "
$ws.Range("D18").Value = "3.171.36"
$ws.Range("E18").Value = "  -8.52%  "
$ws.Range("D19").Value = "'5.73"
$ws.Range("E19").Value = "  -4.17%  "
$ws.Range("D20").Value = "'12.96"
$ws.Range("E20").Value = "  -6.43%  "
$ws.Range("D21").Value = "'352.34"
$ws.Range("E21").Value = "  -5.47%  "
$ws.Range("D22").Value = "'7.17"
$ws.Range("E22").Value = "  -6.81%  "
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("D24").Value = "'68.63"
$ws.Range("E24").Value = "  -6.47%  "
$ws.Range("E25").Value = "  -5.53%  "
$ws.Range("D26").Value = "'0.504"
$ws.Range("E26").Value = "  -6.22%  "
$ws.Range("D27").Value = "'9.58"
$ws.Range("E27").Value = "  -4.89%  "
$ws.Range("E28").Value = "  -0.88%  "
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("E30").Value = "  -0.15%  "
$ws.Range("D31").Value = "'5.53"
$ws.Range("E31").Value = "  -7.48%  "
$ws.Range("D32").Value = "'1.90"
$ws.Range("E32").Value = "  -5.19%  "
$ws.Range("D33").Value = "'21.96"
$ws.Range("E33").Value = "  -7.34%  "
$ws.Range("D34").Value = "'6.64"
$ws.Range("E34").Value = "  -5.84%  "
$ws.Range("E35").Value = "  -5.74%  "
$ws.Range("D36").Value = "'1.44"
$ws.Range("E36").Value = "  -7.75%  "
$ws.Range("D37").Value = "'153.83"
$ws.Range("E37").Value = "  -4.54%  "
$ws.Range("D38").Value = "'0.817"
$ws.Range("E38").Value = "  -7.76%  "
$ws.Range("D39").Value = "'26.02"
$ws.Range("E39").Value = "  -7.73%  "
$ws.Range("E40").Value = "  -6.90%  "
$ws.Range("D41").Value = "'2.48"
$ws.Range("E41").Value = "  -3.35%  "
$ws.Range("D42").Value = "2.608.35"
$ws.Range("E42").Value = "  -6.95%  "
$ws.Range("D43").Value = "'4.18"
$ws.Range("E43").Value = "  -7.57%  "
$ws.Range("D44").Value = "'39.39"
$ws.Range("E44").Value = "  -1.62%  "
$ws.Range("D45").Value = "'5.97"
$ws.Range("E45").Value = "  -8.16%  "
$ws.Range("D46").Value = "'0.0649"
$ws.Range("E46").Value = "  -6.66%  "
$ws.Range("D47").Value = "'23.74"
$ws.Range("E47").Value = "  -5.78%  "
$ws.Range("D48").Value = "'319.19"
$ws.Range("E48").Value = "  -6.26%  "
$ws.Range("D49").Value = "'0.0270"
$ws.Range("E49").Value = "  -7.86%  "
$ws.Range("E50").Value = "  -3.78%  "
$ws.Range("E51").Value = "  +0.00%  "
